$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2-5 down to 3-6
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new data values
$newRow2 = @(1.5685228013090549, 0.46337244269654099, -17.62906530728721, 5.5477245430987079, 458.06222433649594, 9.7108016746029424, 524.33013901089589, 14.32178926258366, 7.2139689645923548, 1.0305669949417651, 0.062753020201732942, 0.062753020201732942, 0.47343571165143011, -0.82425605185365147)

for ($i = 0; $i -lt $newRow2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $newRow2[$i]
}
